$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Directory" column value for the Video row used an inconsistent
# capitalization / legacy directory layout. Fix it to match the actual
# on-disk path used elsewhere in the repository.
$ws.Range("E2").Value = "data/multimedia/video/"

# Restore the cursor/selection position as left by the editor.
$ws.Range("E3").Select()
